# Fill in the newly-drawn lottery numbers for 971회차 (round 971), the
# middle table (columns J:Q) of the third block of rows (15-21), plus the
# "결과"/"보너스" row for the first table (columns A:H) of that same block.
#
# Layout reminder: each block of 7 rows (라인1..5, 결과, 보너스) is repeated
# three times across (A:H, J:Q, S:Z) for three consecutive rounds. Rows
# 16-20 are the five "line" rows and row 21 is the 결과(result)/보너스(bonus)
# row for that block's round.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Middle table (J:Q block, round 971회차) - 5 number lines, rows 16-20
$ws.Range("K16").Value = 6
$ws.Range("L16").Value = 9
$ws.Range("M16").Value = 21
$ws.Range("N16").Value = 28
$ws.Range("O16").Value = 42
$ws.Range("P16").Value = 45

$ws.Range("K17").Value = 8
$ws.Range("L17").Value = 11
$ws.Range("M17").Value = 18
$ws.Range("N17").Value = 28
$ws.Range("O17").Value = 38
$ws.Range("P17").Value = 42

$ws.Range("K18").Value = 17
$ws.Range("L18").Value = 6
$ws.Range("M18").Value = 25
$ws.Range("N18").Value = 30
$ws.Range("O18").Value = 24
$ws.Range("P18").Value = 42

$ws.Range("K19").Value = 18
$ws.Range("L19").Value = 9
$ws.Range("M19").Value = 27
$ws.Range("N19").Value = 21
$ws.Range("O19").Value = 28
$ws.Range("P19").Value = 39

$ws.Range("K20").Value = 12
$ws.Range("L20").Value = 18
$ws.Range("M20").Value = 28
$ws.Range("N20").Value = 32
$ws.Range("O20").Value = 34
$ws.Range("P20").Value = 43

# First table (A:H block, round 970회차) - 결과/보너스 row 21
$ws.Range("B21").Value = 9
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 21
$ws.Range("F21").Value = 28
$ws.Range("G21").Value = 36
$ws.Range("H21").Value = 5

# Move the active selection to match where the user last clicked after
# entering the new numbers.
$ws.Range("P21").Select() | Out-Null
